$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("J1").Value = "number_previous_bleedings"
$ws.Range("J2").Value = 1
$ws.Range("J3").Value = 0
$ws.Range("J4").Value = 0

$ws.Range("J2").Select()
